$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.233197250808075
$ws.Range("D2").Value = 0.1312773035674297
$ws.Range("E2").Value = 1.133037429069134
$ws.Range("F2").Value = 2.883205110301162
$ws.Range("G2").Value = 0.002453524519382204
$ws.Range("I2").Value = 1.005780433012761
$ws.Range("L2").Value = 0.6584819995441649
$ws.Range("M2").Value = 0.4423747884246367

$ws.Range("B3").Value = 1.151319004853804
$ws.Range("D3").Value = 0.1255249957752937
$ws.Range("E3").Value = 0.9882062833839882
$ws.Range("F3").Value = 2.722142274669864
$ws.Range("G3").Value = 0.00246490523351781
$ws.Range("I3").Value = 1.033255439888704
$ws.Range("L3").Value = 0.6038886389110587
$ws.Range("M3").Value = 0.4097149382431198

$ws.Range("B4").Value = 1.101414397856814
$ws.Range("D4").Value = 0.1221805798483331
$ws.Range("E4").Value = 0.8992458334219293
$ws.Range("F4").Value = 2.626493396448012
$ws.Range("G4").Value = 0.002472232856386984
$ws.Range("I4").Value = 1.051068573535993
$ws.Range("L4").Value = 0.5707359580252387
$ws.Range("M4").Value = 0.3898467245697717

$ws.Range("B5").Value = 1.081170176380084
$ws.Range("D5").Value = 0.1208635171718413
$ws.Range("E5").Value = 0.8629775507451711
$ws.Range("F5").Value = 2.588308421535061
$ws.Range("G5").Value = 0.002475304833636589
$ws.Range("I5").Value = 1.05856377288155
$ws.Range("L5").Value = 0.5573157195128147
$ws.Range("M5").Value = 0.3817959705443315

$ws.Range("B6").Value = 1.077814202070414
$ws.Range("D6").Value = 0.1206475482631504
$ws.Range("E6").Value = 0.8569540313695541
$ws.Range("F6").Value = 2.582015047715998
$ws.Range("G6").Value = 0.002475820135499814
$ws.Range("I6").Value = 1.059822581732353
$ws.Range("L6").Value = 0.5550926408157579
$ws.Range("M6").Value = 0.3804618858430331

$ws.Range("B7").Value = 1.101141004003637
$ws.Range("D7").Value = 0.1221626337349164
$ws.Range("E7").Value = 0.8987567822406959
$ws.Range("F7").Value = 2.625975240996382
$ws.Range("G7").Value = 0.002472273937567993
$ws.Range("I7").Value = 1.051168701359979
$ws.Range("L7").Value = 0.5705546081872797
$ws.Range("M7").Value = 0.3897379653016344

$ws.Range("B8").Value = 1.204888720385952
$ws.Range("D8").Value = 0.1292541706045256
$ws.Range("E8").Value = 1.083101837305918
$ws.Range("F8").Value = 2.826983178027888
$ws.Range("G8").Value = 0.002457378340061374
$ws.Range("I8").Value = 1.015057407179969
$ws.Range("L8").Value = 0.6395801590695385
$ws.Range("M8").Value = 0.4310747387634208

$ws.Range("B9").Value = 1.41130149988561
$ws.Range("D9").Value = 0.1447108880791603
$ws.Range("E9").Value = 1.444712242713337
$ws.Range("F9").Value = 3.247985200314844
$ws.Range("G9").Value = 0.002430843171569034
$ws.Range("I9").Value = 0.951773525028619
$ws.Range("L9").Value = 0.7779928598517643
$ws.Range("M9").Value = 0.5136481382919413

$ws.Range("B10").Value = 1.564831988596438
$ws.Range("D10").Value = 0.157101137707258
$ws.Range("E10").Value = 1.711016583969126
$ws.Range("F10").Value = 3.575213028226443
$ws.Range("G10").Value = 0.002412948501376417
$ws.Range("I10").Value = 0.9099309523224512
$ws.Range("L10").Value = 0.8817492119154906
$ws.Range("M10").Value = 0.5753062029799025

$ws.Range("B11").Value = 1.635103367109707
$ws.Range("D11").Value = 0.1629823609494849
$ws.Range("E11").Value = 1.832430983881096
$ws.Range("F11").Value = 3.728320444587013
$ws.Range("G11").Value = 0.00240514875980438
$ws.Range("I11").Value = 0.8919194836049034
$ws.Range("L11").Value = 0.9294454863702981
$ws.Range("M11").Value = 0.6035876625097387

$ws.Range("B12").Value = 1.661776357647682
$ws.Range("D12").Value = 0.1652463098193664
$ws.Range("E12").Value = 1.87845701251527
$ws.Range("F12").Value = 3.786938647137532
$ws.Range("G12").Value = 0.002402243658878866
$ws.Range("I12").Value = 0.8852473259317772
$ws.Range("L12").Value = 0.9475820471472503
$ws.Range("M12").Value = 0.6143318759023089

$ws.Range("B13").Value = 1.656029044763557
$ws.Range("D13").Value = 0.1647570618796408
$ws.Range("E13").Value = 1.868542138915615
$ws.Range("F13").Value = 3.774285250015794
$ws.Range("G13").Value = 0.002402867175613439
$ws.Range("I13").Value = 0.8866776735450692
$ws.Range("L13").Value = 0.9436726198581766
$ws.Range("M13").Value = 0.6120163578441122

$ws.Range("B14").Value = 1.637296506680514
$ws.Range("D14").Value = 0.1631678690800413
$ws.Range("E14").Value = 1.836216531348072
$ws.Range("F14").Value = 3.733130010725944
$ws.Range("G14").Value = 0.002404908786491238
$ws.Range("I14").Value = 0.8913675795785263
$ws.Range("L14").Value = 0.9309360678587666
$ws.Range("M14").Value = 0.6044708944209845

$ws.Range("B15").Value = 1.625830487951987
$ws.Range("D15").Value = 0.1621992923310245
$ws.Range("E15").Value = 1.816422853371876
$ws.Range("F15").Value = 3.708005417217777
$ws.Range("G15").Value = 0.002406165632356921
$ws.Range("I15").Value = 0.8942596469024036
$ws.Range("L15").Value = 0.923144439151514
$ws.Range("M15").Value = 0.5998536262288781

$ws.Range("B16").Value = 1.560248313782722
$ws.Range("D16").Value = 0.1567218526062391
$ws.Range("E16").Value = 1.70308817291675
$ws.Range("F16").Value = 3.565295042686074
$ws.Range("G16").Value = 0.002413465045235582
$ws.Range("I16").Value = 0.911128738372506
$ws.Range("L16").Value = 0.8786424492381855
$ws.Range("M16").Value = 0.5734627263884136

$ws.Range("B17").Value = 1.520126460301583
$ws.Range("D17").Value = 0.1534254878657322
$ws.Range("E17").Value = 1.633637488996044
$ws.Range("F17").Value = 3.478855422149877
$ws.Range("G17").Value = 0.002418029896128213
$ws.Range("I17").Value = 0.9217402898396525
$ws.Range("L17").Value = 0.8514716932942008
$ws.Range("M17").Value = 0.5573332132049842

$ws.Range("B18").Value = 1.497089802025755
$ws.Range("D18").Value = 0.1515524496148686
$ws.Range("E18").Value = 1.593716196198613
$ws.Range("F18").Value = 3.429535888065857
$ws.Range("G18").Value = 0.002420687569351262
$ws.Range("I18").Value = 0.9279399767584646
$ws.Range("L18").Value = 0.8358903514469205
$ws.Range("M18").Value = 0.5480777929687548

$ws.Range("B19").Value = 1.489296891163008
$ws.Range("D19").Value = 0.1509221621459318
$ws.Range("E19").Value = 1.580203510782184
$ws.Range("F19").Value = 3.41290468326747
$ws.Range("G19").Value = 0.002421592938606953
$ws.Range("I19").Value = 0.9300555681574192
$ws.Range("L19").Value = 0.8306226673371384
$ws.Range("M19").Value = 0.5449477852791063

$ws.Range("B20").Value = 1.524393315184454
$ws.Range("D20").Value = 0.153774005152485
$ws.Range("E20").Value = 1.641027984090584
$ws.Range("F20").Value = 3.48801564551519
$ws.Range("G20").Value = 0.002417540642045981
$ws.Range("I20").Value = 0.9206007061039667
$ws.Range("L20").Value = 0.8543592227620422
$ws.Range("M20").Value = 0.559047957983779

$ws.Range("B21").Value = 1.642796996083405
$ws.Range("D21").Value = 0.1636336398364335
$ws.Range("E21").Value = 1.845709933215034
$ws.Range("F21").Value = 3.745200706959139
$ws.Range("G21").Value = 0.002404307803668551
$ws.Range("I21").Value = 0.8899860025842514
$ws.Range("L21").Value = 0.9346750369961967
$ws.Range("M21").Value = 0.6066862304124214

$ws.Range("B22").Value = 1.720546893701055
$ws.Range("D22").Value = 0.1702931151815505
$ws.Range("E22").Value = 1.97977216266122
$ws.Range("F22").Value = 3.917028883939111
$ws.Range("G22").Value = 0.002395941793718285
$ws.Range("I22").Value = 0.8708431080054151
$ws.Range("L22").Value = 0.9876050682991036
$ws.Range("M22").Value = 0.6380230698616742

$ws.Range("B23").Value = 1.679016451007897
$ws.Range("D23").Value = 0.166718539061236
$ws.Range("E23").Value = 1.908190601180422
$ws.Range("F23").Value = 3.824968783743714
$ws.Range("G23").Value = 0.002400381213784904
$ws.Range("I23").Value = 0.8809803923370367
$ws.Range("L23").Value = 0.9593139358477742
$ws.Range("M23").Value = 0.6212790877174683

$ws.Range("B24").Value = 1.522464176939536
$ws.Range("D24").Value = 0.1536163718186003
$ws.Range("E24").Value = 1.637686720736184
$ws.Range("F24").Value = 3.483873140293156
$ws.Range("G24").Value = 0.00241776173026577
$ws.Range("I24").Value = 0.9211156037561103
$ws.Range("L24").Value = 0.8530536479611328
$ws.Range("M24").Value = 0.5582726669740055

$ws.Range("B25").Value = 1.35513695581659
$ws.Range("D25").Value = 0.1403544089792206
$ws.Range("E25").Value = 1.346818441759183
$ws.Range("F25").Value = 3.131062175588255
$ws.Range("G25").Value = 0.002437738342855809
$ws.Range("I25").Value = 0.9680810758473388
$ws.Range("L25").Value = 0.7402010924211311
$ws.Range("M25").Value = 0.4911409602197239
